$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 125.8
$ws.Range("I8").Value = 125.8
$ws.Range("K8").Value = 377.4
$ws.Range("M8").Value = -238.4
$ws.Range("H9").Value = 218.66667
$ws.Range("J9").Value = 212.4
$ws.Range("L9").Value = 212.4
$ws.Range("N9").Value = -550.4
$ws.Range("H51").Value = 11404
$ws.Range("I51").Value = 5373.75
$ws.Range("J51").Value = 14849.857
$ws.Range("K51").Value = 5373.75
$ws.Range("L51").Value = 14849.857
$ws.Range("M51").Value = -4889.75
$ws.Range("N51").Value = -15817.857
$ws.Range("H80").Value = 846.2857
$ws.Range("I80").Value = 660
$ws.Range("K80").Value = 1980
$ws.Range("M80").Value = -982
$ws.Range("H83").Value = 846.2857
$ws.Range("I83").Value = 660
$ws.Range("K83").Value = 5940
$ws.Range("M83").Value = -948
$ws.Range("H98").Value = 1369.6154
$ws.Range("I98").Value = 1122.2
$ws.Range("K98").Value = 1122.2
$ws.Range("M98").Value = 375.8
$ws.Range("H107").Value = 1495.4667
$ws.Range("I107").Value = 1495.4667
$ws.Range("K107").Value = 1495.4667
$ws.Range("M107").Value = 424.5333000000001
$ws.Range("H122").Value = 1369.6154
$ws.Range("I122").Value = 1122.2
$ws.Range("K122").Value = 3366.6
$ws.Range("M122").Value = -916.6000000000004
$ws.Range("H132").Value = 5566.72
$ws.Range("I132").Value = 6469.6313
$ws.Range("K132").Value = 19408.8939
$ws.Range("M132").Value = -16878.8939
$ws.Range("H137").Value = 2205.3438
$ws.Range("J137").Value = 2724.8333
$ws.Range("L137").Value = 8174.499899999999
$ws.Range("N137").Value = -13274.4999
$ws.Range("H138").Value = 3159.3076
$ws.Range("I138").Value = 3096
$ws.Range("J138").Value = 3187.4443
$ws.Range("K138").Value = 9288
$ws.Range("L138").Value = 9562.332900000001
$ws.Range("M138").Value = -4148
$ws.Range("N138").Value = -19842.3329

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11499.667
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 11499.667
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 11499.667
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -11725.667
$ws.Range("H32").Value = 3588.3901
$ws.Range("I32").Value = 3628.1482
$ws.Range("K32").Value = 3628.1482
$ws.Range("M32").Value = -3341.1482
$ws.Range("H61").Value = 9522.037
$ws.Range("I61").Value = 9373.895
$ws.Range("J61").Value = 9873.875
$ws.Range("K61").Value = 9373.895
$ws.Range("L61").Value = 9873.875
$ws.Range("M61").Value = -9161.895
$ws.Range("N61").Value = -10297.875
$ws.Range("H74").Value = 5738.84
$ws.Range("I74").Value = 4233.5884
$ws.Range("K74").Value = 4233.5884
$ws.Range("M74").Value = -3359.5884
$ws.Range("H77").Value = 5738.84
$ws.Range("I77").Value = 4233.5884
$ws.Range("K77").Value = 21167.942
$ws.Range("M77").Value = -16799.942
$ws.Range("H97").Value = 1074.75
$ws.Range("I97").Value = 1182.3
$ws.Range("J97").Value = 537
$ws.Range("K97").Value = 1182.3
$ws.Range("L97").Value = 537
$ws.Range("M97").Value = -686.3
$ws.Range("N97").Value = -1529
$ws.Range("H116").Value = 11499.667
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 11499.667
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 11499.667
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -16087.667
$ws.Range("H132").Value = 3118.3157
$ws.Range("I132").Value = 3163.5151
$ws.Range("K132").Value = 9490.5453
$ws.Range("M132").Value = -6960.5453
$ws.Range("H135").Value = 172000
$ws.Range("J135").Value = 172000
$ws.Range("L135").Value = 172000
$ws.Range("N135").Value = -182140
$ws.Range("H136").Value = 9522.037
$ws.Range("I136").Value = 9373.895
$ws.Range("J136").Value = 9873.875
$ws.Range("K136").Value = 28121.685
$ws.Range("L136").Value = 29621.625
$ws.Range("M136").Value = -25571.685
$ws.Range("N136").Value = -34721.625

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11499.667
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 11499.667
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 11499.667
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -11727.667
$ws.Range("H86").Value = 3516.3333
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 3516.3333
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 4798.375
$ws.Range("I107").Value = 4699
$ws.Range("J107").Value = 5096.5
$ws.Range("K107").Value = 4699
$ws.Range("L107").Value = 5096.5
$ws.Range("M107").Value = -2779
$ws.Range("N107").Value = -8936.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 20000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H31").Value = 5710.125
$ws.Range("I31").Value = 4716.5
$ws.Range("J31").Value = 6419.857
$ws.Range("K31").Value = 4716.5
$ws.Range("L31").Value = 6419.857
$ws.Range("M31").Value = -4421.5
$ws.Range("N31").Value = -7009.857
$ws.Range("H34").Value = 5710.125
$ws.Range("I34").Value = 4716.5
$ws.Range("J34").Value = 6419.857
$ws.Range("K34").Value = 4716.5
$ws.Range("L34").Value = 6419.857
$ws.Range("M34").Value = -4514.5
$ws.Range("N34").Value = -6823.857
$ws.Range("H41").Value = 22999.6
$ws.Range("J41").Value = 23666.666
$ws.Range("L41").Value = 23666.666
$ws.Range("N41").Value = -24522.666
$ws.Range("H59").Value = 40700
$ws.Range("J59").Value = 40700
$ws.Range("L59").Value = 40700
$ws.Range("N59").Value = -42990
$ws.Range("H60").Value = 15612.5
$ws.Range("J60").Value = 28500
$ws.Range("L60").Value = 28500
$ws.Range("N60").Value = -29522
$ws.Range("H74").Value = 42527
$ws.Range("J74").Value = 42527
$ws.Range("L74").Value = 42527
$ws.Range("N74").Value = -44275
$ws.Range("H77").Value = 42527
$ws.Range("J77").Value = 42527
$ws.Range("L77").Value = 127581
$ws.Range("N77").Value = -136317
$ws.Range("H138").Value = 133460.31
$ws.Range("I138").Value = 84999
$ws.Range("J138").Value = 399997.5
$ws.Range("K138").Value = 84999
$ws.Range("L138").Value = 399997.5
$ws.Range("M138").Value = -79859
$ws.Range("N138").Value = -410277.5
$ws.Range("H141").Value = 38115.25
$ws.Range("I141").Value = 37887
$ws.Range("J141").Value = 38800
$ws.Range("K141").Value = 37887
$ws.Range("L141").Value = 38800
$ws.Range("M141").Value = -32707
$ws.Range("N141").Value = -49160

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4750
$ws.Range("I5").Value = 7000
$ws.Range("K5").Value = 21000
$ws.Range("M5").Value = -20888
$ws.Range("H81").Value = 1029.5
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 1029.5
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 3088.5
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -5334.5
$ws.Range("H84").Value = 1029.5
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 1029.5
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 9265.5
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -20497.5
$ws.Range("H122").Value = 544.6
$ws.Range("I122").Value = 328.66666
$ws.Range("J122").Value = 688.55554
$ws.Range("K122").Value = 2957.99994
$ws.Range("L122").Value = 6196.99986
$ws.Range("M122").Value = -507.9999399999997
$ws.Range("N122").Value = -11096.99986
$ws.Range("H128").Value = 1471694.8
$ws.Range("I128").Value = 1471694.8
$ws.Range("K128").Value = 4415084.4
$ws.Range("M128").Value = -4410104.4
$ws.Range("H135").Value = 4750
$ws.Range("I135").Value = 7000
$ws.Range("K135").Value = 63000
$ws.Range("M135").Value = -60465
$ws.Range("H139").Value = 3520.2778
$ws.Range("I139").Value = 3633.8572
$ws.Range("J139").Value = 3122.75
$ws.Range("K139").Value = 10901.5716
$ws.Range("L139").Value = 9368.25
$ws.Range("M139").Value = -5761.571599999999
$ws.Range("N139").Value = -19648.25
$ws.Range("H140").Value = 1368.5
$ws.Range("I140").Value = 1108.72
$ws.Range("K140").Value = 3326.16
$ws.Range("M140").Value = 1853.84

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 51614
$ws.Range("J46").Value = 51614
$ws.Range("L46").Value = 51614
$ws.Range("N46").Value = -51926

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1636.909
$ws.Range("J22").Value = 1337.8334
$ws.Range("L22").Value = 1337.8334
$ws.Range("N22").Value = -1927.8334
$ws.Range("H27").Value = 1636.909
$ws.Range("J27").Value = 1337.8334
$ws.Range("L27").Value = 1337.8334
$ws.Range("N27").Value = -1551.8334
$ws.Range("H132").Value = 11125.191
$ws.Range("J132").Value = 8706.462
$ws.Range("L132").Value = 26119.386
$ws.Range("N132").Value = -31179.386

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 18169.8
$ws.Range("I41").Value = 13000
$ws.Range("J41").Value = 18744.223
$ws.Range("K41").Value = 13000
$ws.Range("L41").Value = 18744.223
$ws.Range("M41").Value = -12610
$ws.Range("N41").Value = -19524.223
$ws.Range("H62").Value = 5741.625
$ws.Range("I62").Value = 5299.6924
$ws.Range("K62").Value = 5299.6924
$ws.Range("M62").Value = -4675.6924
$ws.Range("H65").Value = 5741.625
$ws.Range("I65").Value = 5299.6924
$ws.Range("K65").Value = 26498.462
$ws.Range("M65").Value = -23378.462
